# Add a new "2022-Q3" quarterly sheet ahead of the existing "2022-Q2" sheet,
# populate it with the new quarter's fund-position data, and update the
# "总计" (summary) sheet with the corresponding new row.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (same
#        layout/headers/styles) and inserting it immediately before it. ---
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q3"

# --- 2. Fill in the 2022-Q3 figures (fund code/name stay the same; the
#        scale/position/ranking numbers are the new quarter's values). ---
# Row 2: fund 008763 (A share)
$newSheet.Range("D2").Value = "'20.44"
$newSheet.Range("E2").Value = "'90.19"
$newSheet.Range("F2").Value = "'4.86"
$newSheet.Range("G2").Value = "'0.9934"
$newSheet.Range("H2").Value = 10

# Row 3: fund 008764 (C share)
$newSheet.Range("D3").Value = "'15.02"
$newSheet.Range("E3").Value = "'90.19"
$newSheet.Range("F3").Value = "'4.86"
$newSheet.Range("G3").Value = "'0.7300"
$newSheet.Range("H3").Value = 10

# --- 3. Update the "总计" summary sheet: shift the existing quarters down
#        one row and insert the new 2022-Q3 total at the top (row 2). The
#        previously-last row (2021-Q2) is re-appended at the new bottom row
#        so no historical quarter is lost. ---
$summary = $wb.Worksheets.Item("总计")

# Duplicate the last existing row's formatting into the new row 7 first.
$summary.Range("A6:D6").Copy($summary.Range("A7:D7"))
$summary.Range("A7").Value = 5

# Shift rows 6->5->4->3->2 down by one (oldest first) so each row keeps the
# label/value that belongs one quarter further back, freeing up row 2 for
# the brand-new quarter. (Read with Value2 -- plain .Value round-trips as an
# opaque wrapper in this host when re-assigned straight into another cell.)
$summary.Range("B7").Value = $summary.Range("B6").Value2
$summary.Range("D7").Value = $summary.Range("D6").Value2

$summary.Range("B6").Value = $summary.Range("B5").Value2
$summary.Range("D6").Value = $summary.Range("D5").Value2

$summary.Range("B5").Value = $summary.Range("B4").Value2
$summary.Range("D5").Value = $summary.Range("D4").Value2

$summary.Range("B4").Value = $summary.Range("B3").Value2
$summary.Range("D4").Value = $summary.Range("D3").Value2

$summary.Range("B3").Value = $summary.Range("B2").Value2
$summary.Range("D3").Value = $summary.Range("D2").Value2

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 1.72

# --- 4. Restore the originally-selected tab (2021-Q2) as the active sheet. ---
$wb.Worksheets.Item("2021-Q2").Activate()
